{"js": "// Insert three new paragraphs (blank, \"COURT LOCATION: ...\" , blank)\n// immediately after the \"DATE:   The _______ day of ...\" paragraph and\n// before the \"PLAINTIFF/TENANT\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the paragraph that contains the DATE placeholder text. Searching\n// by content (rather than a hard-coded index) keeps the script resilient\n// to any unrelated changes earlier in the document.\nlet dateParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"DATE:\") !== -1 && text.indexOf(\"o'clock\") !== -1) {\n    dateParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!dateParagraph) {\n  throw new Error('Could not find the \"DATE:\" paragraph to insert after.');\n}\n\n// Insert the new blank paragraph right after the DATE paragraph. It\n// inherits the DATE paragraph's paragraph formatting (tab stops, Garamond\n// font), matching the target markup.\nconst blankBefore = dateParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\n// Insert the \"COURT LOCATION\" paragraph after that blank paragraph.\nconst courtLocationParagraph = blankBefore.insertParagraph(\n  \"COURT LOCATION: {{ trial_court.address.on_one_line() }}\",\n  Word.InsertLocation.after\n);\n\n// Insert a trailing blank paragraph after the COURT LOCATION paragraph.\ncourtLocationParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert three new paragraphs (blank, \"COURT LOCATION: ...\" , blank)\n# immediately after the \"DATE:   The _______ day of ...\" paragraph and\n# before the \"PLAINTIFF/TENANT\" paragraph.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"DATE:*o'clock*\", $false, $false, $true)\n\nif (-not $found) {\n    throw \"Could not find the 'DATE:' paragraph to insert after.\"\n}\n\n$dateParagraph = $rng.Paragraphs(1)\n\n# Insert a blank paragraph right after the DATE paragraph. It inherits the\n# DATE paragraph's formatting (tab stops, Garamond font), matching the\n# target markup.\n$dateParagraph.Range.InsertParagraphAfter()\n\n# Insert the \"COURT LOCATION\" paragraph after that blank paragraph.\n$blankBefore = $dateParagraph.Next()\n$blankBefore.Range.InsertParagraphAfter()\n\n$courtLocationParagraph = $blankBefore.Next()\n$courtLocationParagraph.Range.Text = \"COURT LOCATION: {{ trial_court.address.on_one_line() }}\"\n\n# Insert a trailing blank paragraph after the COURT LOCATION paragraph.\n$courtLocationParagraph.Range.InsertParagraphAfter()\n"}
